# Fixed error in sinuous length calculation and summary calculations.
# Fixed a small bug in if statements for creek ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (RS order 1 -> 3, values recalculated)
$ws.Range("A2").Value = 3
$ws.Range("C2").Value = 250.1787155501902
$ws.Range("D2").Value = 10.87733545870392
$ws.Range("E2").Value = 1.053701939707178
$ws.Range("F2").Value = 5.019123444325006
$ws.Range("G2").Value = 0.40497819820176
$ws.Range("H2").Value = 0.6252391042916673
$ws.Range("I2").Value = 234.276888520289
$ws.Range("J2").Value = 10.18595167479518

# Row 3 (RS order 2, values recalculated; A3/B3 unchanged)
$ws.Range("C3").Value = 446.4751801064722
$ws.Range("D3").Value = 148.8250600354907
$ws.Range("E3").Value = 1.136229726332538

# Row 4 (RS order 3 -> 1, values recalculated)
$ws.Range("A4").Value = 1
$ws.Range("C4").Value = 370.6883835420692
$ws.Range("D4").Value = 370.6883835420692
$ws.Range("E4").Value = 1.679318423721673
